$wb = $excel.ActiveWorkbook

# The same updates need to be applied to both the "展览" sheet and the
# "全部类型" sheet, since they contain duplicated data (F column = "想去人数").
$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    3  = 96
    4  = 393
    5  = 11574
    6  = 795
    8  = 15
    11 = 167
    12 = 21
    18 = 1328
    19 = 76
    20 = 901
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
